# Updates cryptos list values (Price / Volume(1h)) per the scraped-data refresh,
# and swaps the RenderToken/Bittensor row order (rows 40-41).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "0.0498", "1.00")
    # are kept verbatim instead of being coerced into numbers, then restore
    # the default (unstyled) cell style so no formatting is introduced.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText "D2" "58.846.63"
Set-CellText "E2" "  +0.92%  "
Set-CellText "D3" "2.505.38"
Set-CellText "E3" "  +2.27%  "
Set-CellText "D4" "0.999"
Set-CellText "E4" "  -0.09%  "
Set-CellText "D5" "533.07"
Set-CellText "E5" "  +0.96%  "
Set-CellText "D6" "135.86"
Set-CellText "E6" "  +1.72%  "
Set-CellText "D7" "0.999"
Set-CellText "E7" "  -0.03%  "
Set-CellText "E8" "  +1.60%  "
Set-CellText "E9" "  +2.20%  "
Set-CellText "E10" "  -1.51%  "
Set-CellText "D11" "5.39"
Set-CellText "E11" "  +1.90%  "
Set-CellText "E12" "  +1.44%  "
Set-CellText "D13" "2.941.23"
Set-CellText "E13" "  +1.90%  "
Set-CellText "D14" "58.785.91"
Set-CellText "E14" "  +0.93%  "
Set-CellText "D15" "22.63"
Set-CellText "E15" "  +0.14%  "
Set-CellText "E16" "  +0.12%  "
Set-CellText "D17" "2.500.20"
Set-CellText "E17" "  +1.61%  "
Set-CellText "D18" "10.99"
Set-CellText "E18" "  +2.29%  "
Set-CellText "E19" "  +1.03%  "
Set-CellText "D20" "321.83"
Set-CellText "E20" "  +0.45%  "
Set-CellText "E21" "  +0.27%  "
Set-CellText "E22" "  +4.38%  "
Set-CellText "E23" "  +4.29%  "
Set-CellText "E24" "  +2.58%  "
Set-CellText "E25" "  +0.39%  "
Set-CellText "E26" "  +1.42%  "
Set-CellText "E27" "  +0.68%  "
Set-CellText "E28" "  +1.13%  "
Set-CellText "D29" "171.18"
Set-CellText "E29" "  +4.33%  "
Set-CellText "E30" "  -0.93%  "
Set-CellText "E31" "  -0.39%  "
Set-CellText "E32" "  +4.62%  "
Set-CellText "E33" "  -0.03%  "
Set-CellText "D34" "18.31"
Set-CellText "E34" "  +0.65%  "
Set-CellText "D35" "1.34"
Set-CellText "E35" "  -0.38%  "
Set-CellText "E36" "  -0.10%  "
Set-CellText "E37" "  -0.40%  "
Set-CellText "E38" "  +0.49%  "
Set-CellText "D39" "0.789"
Set-CellText "E39" "  -1.69%  "
Set-CellText "B40" "Bittensor"
Set-CellText "C40" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText "D40" "280.24"
Set-CellText "E40" "  +2.40%  "
Set-CellText "B41" "RenderToken"
Set-CellText "C41" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D41" "5.08"
Set-CellText "E41" "  +0.41%  "
Set-CellText "D42" "0.999"
Set-CellText "E42" "  +0.12%  "
Set-CellText "E43" "  +3.26%  "
Set-CellText "D44" "129.49"
Set-CellText "E44" "  +7.29%  "
Set-CellText "D45" "10.91"
Set-CellText "E45" "  +0.69%  "
Set-CellText "E46" "  -0.01%  "
Set-CellText "D47" "0.0498"
Set-CellText "E47" "  -1.05%  "
Set-CellText "E48" "  +0.19%  "
Set-CellText "D49" "17.14"
Set-CellText "E49" "  +0.57%  "
Set-CellText "D50" "1.750.12"
Set-CellText "E50" "  +0.56%  "
Set-CellText "E51" "  +0.14%  "
